$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column A is treated as text (inline/shared string) for all data rows
$ws.Range("A2:A17").NumberFormat = "@"

$ws.Cells.Item(2, 1).Value = "0"
$ws.Cells.Item(2, 2).Value = 23.66698125826089
$ws.Cells.Item(2, 3).Value = 23.42013439223523
$ws.Cells.Item(2, 4).Value = 0.5628382013489843
$ws.Cells.Item(2, 5).Value = 0.5804712546134158

$ws.Cells.Item(3, 1).Value = "1"
$ws.Cells.Item(3, 2).Value = 23.6652446485445
$ws.Cells.Item(3, 3).Value = 24.05579222437616
$ws.Cells.Item(3, 4).Value = 0.5628585081383963
$ws.Cells.Item(3, 5).Value = 0.613742543171837

$ws.Cells.Item(4, 1).Value = "2"
$ws.Cells.Item(4, 2).Value = 23.66542198747178
$ws.Cells.Item(4, 3).Value = 25.43849144375257
$ws.Cells.Item(4, 4).Value = 0.5615749372769705
$ws.Cells.Item(4, 5).Value = 0.6842200997350938

$ws.Cells.Item(5, 1).Value = "3"
$ws.Cells.Item(5, 2).Value = 23.6659291784576
$ws.Cells.Item(5, 3).Value = 29.09421076393417
$ws.Cells.Item(5, 4).Value = 0.5617109537252648
$ws.Cells.Item(5, 5).Value = 0.8316405635363789

$ws.Cells.Item(6, 1).Value = "4"
$ws.Cells.Item(6, 2).Value = 23.61109032665684
$ws.Cells.Item(6, 3).Value = 22.7552839708626
$ws.Cells.Item(6, 4).Value = 0.5630425821079695
$ws.Cells.Item(6, 5).Value = 0.5316978520688432

$ws.Cells.Item(7, 1).Value = "5"
$ws.Cells.Item(7, 2).Value = 23.56698259185531
$ws.Cells.Item(7, 3).Value = 23.31498834013616
$ws.Cells.Item(7, 4).Value = 0.560701678971915
$ws.Cells.Item(7, 5).Value = 0.565062299647228

$ws.Cells.Item(8, 1).Value = "6"
$ws.Cells.Item(8, 2).Value = 23.61746968066501
$ws.Cells.Item(8, 3).Value = 24.7460004079343
$ws.Cells.Item(8, 4).Value = 0.5623654053350033
$ws.Cells.Item(8, 5).Value = 0.6464005945857576

$ws.Cells.Item(9, 1).Value = "7"
$ws.Cells.Item(9, 2).Value = 23.67097614796583
$ws.Cells.Item(9, 3).Value = 28.35847660851668
$ws.Cells.Item(9, 4).Value = 0.5623341805466151
$ws.Cells.Item(9, 5).Value = 0.8031275800649054

$ws.Cells.Item(10, 1).Value = "9"
$ws.Cells.Item(10, 2).Value = 23.58185794894158
$ws.Cells.Item(10, 3).Value = 28.29419717316968
$ws.Cells.Item(10, 4).Value = 0.5592424193153306
$ws.Cells.Item(10, 5).Value = 0.8042189327551276

$ws.Cells.Item(11, 1).Value = "10"
$ws.Cells.Item(11, 2).Value = 23.53296335514261
$ws.Cells.Item(11, 3).Value = 24.70011184212721
$ws.Cells.Item(11, 4).Value = 0.5607680889382359
$ws.Cells.Item(11, 5).Value = 0.6485018772956437

$ws.Cells.Item(12, 1).Value = "11"
$ws.Cells.Item(12, 2).Value = 23.53707245832457
$ws.Cells.Item(12, 3).Value = 23.336184782947
$ws.Cells.Item(12, 4).Value = 0.5618308869854264
$ws.Cells.Item(12, 5).Value = 0.5707600925522629

$ws.Cells.Item(13, 1).Value = "12"
$ws.Cells.Item(13, 2).Value = 23.62904466930573
$ws.Cells.Item(13, 3).Value = 22.87548817301593
$ws.Cells.Item(13, 4).Value = 0.563244054210409
$ws.Cells.Item(13, 5).Value = 0.5375440827307707

$ws.Cells.Item(14, 1).Value = "13"
$ws.Cells.Item(14, 2).Value = 23.60078928168511
$ws.Cells.Item(14, 3).Value = 29.09656297293069
$ws.Cells.Item(14, 4).Value = 0.5629599752400475
$ws.Cells.Item(14, 5).Value = 0.8334006987874809

$ws.Cells.Item(15, 1).Value = "14"
$ws.Cells.Item(15, 2).Value = 23.59847117915199
$ws.Cells.Item(15, 3).Value = 25.4901777388197
$ws.Cells.Item(15, 4).Value = 0.5619946590095188
$ws.Cells.Item(15, 5).Value = 0.6862067042922995

$ws.Cells.Item(16, 1).Value = "15"
$ws.Cells.Item(16, 2).Value = 23.59772295605822
$ws.Cells.Item(16, 3).Value = 24.04720191284321
$ws.Cells.Item(16, 4).Value = 0.5586627486975914
$ws.Cells.Item(16, 5).Value = 0.6104305880609839

$ws.Cells.Item(17, 1).Value = "16"
$ws.Cells.Item(17, 2).Value = 23.59866594249343
$ws.Cells.Item(17, 3).Value = 23.37781593296855
$ws.Cells.Item(17, 4).Value = 0.5600088516890556
$ws.Cells.Item(17, 5).Value = 0.5763392919183502

Write-Host "Done"
